$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Version" column before column A, shifting Code/Description/Definition
# one column to the right (B/C/D).
$ws.Columns.Item(1).Insert()

# New header cell
$ws.Range("A1").Value = "Version"

# Build the "1.0" text value in a scratch cell via a formula (so Excel treats it
# as text, not a number) and paste it as a value into every data row of the new
# column. This avoids forcing a text NumberFormat onto the destination cells.
$helper = $ws.Range("Z1")
$helper.Formula = "=""1.0"""
$helper.Copy()
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)
}
$helper.Clear()
$excel.CutCopyMode = $false
